$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(0.9519518233928137, 0.6608182499140475, 1.2337267539856231, 3.1575558330682485)"
$ws.Range("C2").Value = "JSU(-0.931632726065527, 1.0241196670198187, 4.356308690931501, 4.460613376122927)"
$ws.Range("D2").Value = "JSU(-0.7230186148342803, 0.9887969211400421, 0.9308941180233206, 2.2503016568178387)"
$ws.Range("E2").Value = "NIG(1.0143194977821497, 0.7852580435735375, 4.049792238410902, 5.693886138403583)"
